$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1742
$ws.Range("F4").Value = 160
$ws.Range("F5").Value = 436
$ws.Range("F6").Value = 811
$ws.Range("F7").Value = 242
$ws.Range("F8").Value = 1162
$ws.Range("F9").Value = 332
$ws.Range("F11").Value = 872
$ws.Range("F12").Value = 674
$ws.Range("F18").Value = 2900
$ws.Range("F19").Value = 2612
$ws.Range("F23").Value = 314
$ws.Range("F31").Value = 297
$ws.Range("F32").Value = 1081

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1112
$ws.Range("F10").Value = 26
$ws.Range("F26").Value = 3901
$ws.Range("F33").Value = 159

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1777
$ws.Range("F5").Value = 2439
$ws.Range("F9").Value = 1305
$ws.Range("F10").Value = 355

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1777
$ws.Range("F4").Value = 2439
$ws.Range("F5").Value = 1742
$ws.Range("F7").Value = 1305
$ws.Range("F8").Value = 355
$ws.Range("F10").Value = 160
$ws.Range("F11").Value = 436
$ws.Range("F12").Value = 811
$ws.Range("F13").Value = 242
$ws.Range("F14").Value = 1162
$ws.Range("F15").Value = 332
$ws.Range("F16").Value = 872
$ws.Range("F17").Value = 674
$ws.Range("F18").Value = 1112
$ws.Range("F19").Value = 1112
$ws.Range("F23").Value = 2900
$ws.Range("F24").Value = 2612
$ws.Range("F27").Value = 314
$ws.Range("F28").Value = 26
$ws.Range("F38").Value = 297
$ws.Range("F47").Value = 1081
$ws.Range("F49").Value = 159
